$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1975308641975309
$ws.Range("C2").Value = 0.5462962962962963
$ws.Range("J2").Value = 0.0308641975308642
$ws.Range("P2").Value = 0.1234567901234568
$ws.Range("S2").Value = 0.1018518518518518
# Row 3
$ws.Range("B3").Value = 0.01630434782608696
$ws.Range("C3").Value = 0.02173913043478261
$ws.Range("J3").Value = 0.02173913043478261
$ws.Range("P3").Value = 0.6902173913043478
$ws.Range("S3").Value = 0.25
# Row 4
$ws.Range("J4").Value = 0.1282051282051282
$ws.Range("P4").Value = 0.5128205128205128
$ws.Range("S4").Value = 0.358974358974359
# Row 5
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.3333333333333333
# Row 6
$ws.Range("B6").Value = 0.07818930041152264
$ws.Range("D6").Value = 0.01234567901234568
$ws.Range("F6").Value = 0.02880658436213992
$ws.Range("J6").Value = 0.2345679012345679
$ws.Range("O6").Value = 0.01646090534979424
$ws.Range("Q6").Value = 0.205761316872428
$ws.Range("R6").Value = 0.06584362139917696
$ws.Range("S6").Value = 0.3580246913580247
# Row 7
$ws.Range("B7").Value = 0.1341991341991342
$ws.Range("D7").Value = 0.008658008658008658
$ws.Range("E7").Value = 0.008658008658008658
$ws.Range("F7").Value = 0.0735930735930736
$ws.Range("J7").Value = 0.1385281385281385
$ws.Range("O7").Value = 0.01731601731601732
$ws.Range("Q7").Value = 0.1731601731601732
$ws.Range("R7").Value = 0.07792207792207792
$ws.Range("S7").Value = 0.3679653679653679
# Row 8
$ws.Range("B8").Value = 0.09387755102040816
$ws.Range("D8").Value = 0.00816326530612245
$ws.Range("F8").Value = 0.07346938775510205
$ws.Range("J8").Value = 0.1204081632653061
$ws.Range("O8").Value = 0.02040816326530612
$ws.Range("Q8").Value = 0.1530612244897959
$ws.Range("R8").Value = 0.08163265306122448
$ws.Range("S8").Value = 0.4489795918367347
# Row 9
$ws.Range("B9").Value = 0.05913978494623656
$ws.Range("D9").Value = 0.01075268817204301
$ws.Range("F9").Value = 0.08602150537634409
$ws.Range("J9").Value = 0.08064516129032258
$ws.Range("O9").Value = 0.01612903225806452
$ws.Range("Q9").Value = 0.2419354838709677
$ws.Range("R9").Value = 0.07526881720430108
$ws.Range("S9").Value = 0.4301075268817204
# Row 10
$ws.Range("B10").Value = 0.1029207232267038
$ws.Range("D10").Value = 0.02086230876216968
$ws.Range("E10").Value = 0.0006954102920723226
$ws.Range("F10").Value = 0.06954102920723226
$ws.Range("J10").Value = 0.1230876216968011
$ws.Range("O10").Value = 0.01460361613351878
$ws.Range("Q10").Value = 0.1974965229485396
$ws.Range("R10").Value = 0.08762169680111266
$ws.Range("S10").Value = 0.3831710709318498
# Row 11
$ws.Range("G11").Value = 0.1640625
$ws.Range("J11").Value = 0.08072916666666667
$ws.Range("K11").Value = 0.2057291666666667
$ws.Range("L11").Value = 0.5364583333333334
$ws.Range("S11").Value = 0.01302083333333333
# Row 12
$ws.Range("G12").Value = 0.6698564593301436
$ws.Range("J12").Value = 0.2583732057416268
$ws.Range("K12").Value = 0.004784688995215311
$ws.Range("L12").Value = 0.03349282296650718
$ws.Range("S12").Value = 0.03349282296650718
# Row 13
$ws.Range("G13").Value = 0.7560975609756098
$ws.Range("J13").Value = 0.1463414634146341
$ws.Range("S13").Value = 0.0975609756097561
# Row 15
$ws.Range("F15").Value = 0.03381642512077294
$ws.Range("H15").Value = 0.178743961352657
$ws.Range("I15").Value = 0.06280193236714976
$ws.Range("J15").Value = 0.3381642512077295
$ws.Range("K15").Value = 0.0966183574879227
$ws.Range("M15").Value = 0.02415458937198068
$ws.Range("N15").Value = 0.004830917874396135
$ws.Range("O15").Value = 0.03381642512077294
$ws.Range("S15").Value = 0.2270531400966184
# Row 16
$ws.Range("F16").Value = 0.01630434782608696
$ws.Range("H16").Value = 0.1847826086956522
$ws.Range("I16").Value = 0.04891304347826087
$ws.Range("J16").Value = 0.4130434782608696
$ws.Range("K16").Value = 0.1304347826086956
$ws.Range("O16").Value = 0.03260869565217391
$ws.Range("S16").Value = 0.1739130434782609
# Row 17
$ws.Range("F17").Value = 0.01855670103092784
$ws.Range("H17").Value = 0.1731958762886598
$ws.Range("I17").Value = 0.08865979381443299
$ws.Range("J17").Value = 0.4309278350515464
$ws.Range("K17").Value = 0.1195876288659794
$ws.Range("M17").Value = 0.01237113402061856
$ws.Range("N17").Value = 0.002061855670103093
$ws.Range("O17").Value = 0.04329896907216495
$ws.Range("S17").Value = 0.111340206185567
# Row 18
$ws.Range("F18").Value = 0.01401869158878505
$ws.Range("H18").Value = 0.1775700934579439
$ws.Range("I18").Value = 0.1168224299065421
$ws.Range("J18").Value = 0.4532710280373832
$ws.Range("K18").Value = 0.07009345794392523
$ws.Range("M18").Value = 0.004672897196261682
$ws.Range("O18").Value = 0.06074766355140187
$ws.Range("S18").Value = 0.102803738317757
# Row 19
$ws.Range("F19").Value = 0.01101169993117687
$ws.Range("H19").Value = 0.2050929112181693
$ws.Range("I19").Value = 0.06607019958706126
$ws.Range("J19").Value = 0.3840330350997935
$ws.Range("K19").Value = 0.1266345492085341
$ws.Range("M19").Value = 0.01995870612525809
$ws.Range("O19").Value = 0.06194081211286993
$ws.Range("S19").Value = 0.125258086717137
